$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 245; A = 44319; B = 0; C = 1; D = 83.40283569641367 },
    @{ Row = 246; A = 44320; B = 0; C = 1; D = 83.40283569641367 },
    @{ Row = 247; A = 44321; B = 0; C = 1; D = 83.40283569641367 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item(244, 1).Copy()
    $ws.Cells.Item($r.Row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
